$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "d"

$ws.Range("B4").Select()
